$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 180 (the "totals" block
# shifts down from 180-183 to 182-185, and two fresh data rows appear at 180-181).
$ws.Range("A180:G180").EntireRow.Insert()
$ws.Range("A180:G180").EntireRow.Insert()

# New row 180: a complete entry (2014-08-11, 11:15 -> 12:15)
$ws.Range("A180").Value = 2014
$ws.Range("B180").Value = 8
$ws.Range("C180").Value = 11
$ws.Range("D180").Value = 0.46875
$ws.Range("E180").Value = 0.51041666666666663

# Extend the shared "minutes spent" / "hours spent" formulas down through the
# new row 180 (they previously covered F176:F179 / G176:G179).
$ws.Range("F176:F180").Formula = "=(E176-D176)*24*60"
$ws.Range("G176:G180").Formula = "=F176/60"

# New row 181: a started-but-not-finished entry (2014-08-11, 16:00 -> ?)
$ws.Range("A181").Value = 2014
$ws.Range("B181").Value = 8
$ws.Range("C181").Value = 11
$ws.Range("D181").Value = 0.66666666666666663
# No end time yet, so no elapsed-time formula for this row.
$ws.Range("G181").Clear()

$ws.Range("E181").Select()
